$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "GHIJK464748"
$ws.Range("B2").Value = "male"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 86
$ws.Range("G2").Value = 75
$ws.Range("H2").Value = 95
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 78
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 5
$ws.Range("AD2").Value = 9
$ws.Range("AI2").Value = 9
$ws.Range("AN2").Value = 19
$ws.Range("AO2").Value = 41
$ws.Range("AP2").Value = 16
$ws.Range("AQ2").Value = 76
$ws.Range("AR2").Value = $True
$ws.Range("AS2").Value = 2

# Row 3
$ws.Range("A3").Value = "LMNOPQ606162"
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = 75
$ws.Range("H3").Value = 77
$ws.Range("I3").Value = 61
$ws.Range("O3").Value = 10
$ws.Range("Q3").Value = 0
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 7
$ws.Range("AJ3").Value = 9
$ws.Range("AK3").ClearContents()
$ws.Range("AN3").Value = 18
$ws.Range("AO3").Value = 39
$ws.Range("AP3").Value = 3
$ws.Range("AQ3").Value = 60
$ws.Range("AR3").Value = $False
$ws.Range("AS3").Value = 8

# Row 4
$ws.Range("A4").Value = "DEFGHI656667"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 59
$ws.Range("G4").Value = 74
$ws.Range("H4").Value = 55
$ws.Range("J4").Value = 67
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 10
$ws.Range("Q4").ClearContents()
$ws.Range("R4").Value = 7
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 5
$ws.Range("AI4").ClearContents()
$ws.Range("AJ4").Value = 10
$ws.Range("AN4").Value = 19
$ws.Range("AO4").Value = 31
$ws.Range("AP4").Value = 22
$ws.Range("AQ4").Value = 72
$ws.Range("AS4").Value = 3

# Row 5
$ws.Range("A5").Value = "KLDEFG434445"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 78
$ws.Range("H5").Value = 65
$ws.Range("K5").ClearContents()
$ws.Range("L5").Value = 77
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 10
$ws.Range("AC5").Value = 10
$ws.Range("AI5").ClearContents()
$ws.Range("AK5").Value = 10
$ws.Range("AN5").Value = 18
$ws.Range("AP5").Value = 27
$ws.Range("AQ5").Value = 79
$ws.Range("AS5").Value = 1

# Row 6
$ws.Range("A6").Value = "DEFGHI838485"
$ws.Range("E6").Value = 25
$ws.Range("F6").Value = 71
$ws.Range("G6").Value = 56
$ws.Range("H6").Value = 55
$ws.Range("I6").ClearContents()
$ws.Range("L6").Value = 94
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 6
$ws.Range("Q6").ClearContents()
$ws.Range("R6").Value = 10
$ws.Range("AD6").Value = 7
$ws.Range("AI6").Value = 8
$ws.Range("AJ6").ClearContents()
$ws.Range("AN6").Value = 16
$ws.Range("AO6").Value = 34
$ws.Range("AP6").Value = 14
$ws.Range("AQ6").Value = 64
$ws.Range("AR6").Value = $False
$ws.Range("AS6").Value = 6

# Row 7
$ws.Range("A7").Value = "HIJKLM757677"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 79
$ws.Range("G7").Value = 67
$ws.Range("L7").Value = 54
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 10
$ws.Range("Q7").Value = 5
$ws.Range("AC7").Value = 0
$ws.Range("AD7").Value = 6
$ws.Range("AI7").ClearContents()
$ws.Range("AJ7").Value = 0
$ws.Range("AN7").Value = 12
$ws.Range("AO7").Value = 34
$ws.Range("AQ7").Value = 64
$ws.Range("AR7").Value = $False
$ws.Range("AS7").Value = 6

# Row 8
$ws.Range("A8").Value = "DLMNOP808182"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 21
$ws.Range("F8").Value = 85
$ws.Range("G8").Value = 60
$ws.Range("H8").Value = 89
$ws.Range("K8").ClearContents()
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 7
$ws.Range("Q8").Value = 7
$ws.Range("R8").ClearContents()
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 8
$ws.Range("AI8").ClearContents()
$ws.Range("AJ8").Value = 10
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 37
$ws.Range("AP8").Value = 11
$ws.Range("AQ8").Value = 69
$ws.Range("AS8").Value = 5

# Row 9
$ws.Range("A9").Value = "FGHI121314"
$ws.Range("E9").Value = 25
$ws.Range("F9").Value = 94
$ws.Range("G9").Value = 52
$ws.Range("H9").Value = 93
$ws.Range("L9").Value = 87
$ws.Range("N9").Value = 7
$ws.Range("O9").Value = 8
$ws.Range("Q9").Value = 7
$ws.Range("AC9").Value = 7
$ws.Range("AD9").Value = 7
$ws.Range("AI9").Value = 10
$ws.Range("AJ9").ClearContents()
$ws.Range("AN9").Value = 18
$ws.Range("AO9").Value = 40
$ws.Range("AP9").Value = 12
$ws.Range("AQ9").Value = 70
$ws.Range("AS9").Value = 4

# Row 10
$ws.Range("A10").Value = "NOP959697"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = 87
$ws.Range("G10").Value = 59
$ws.Range("H10").Value = 99
$ws.Range("M10").Value = 47
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 8
$ws.Range("R10").Value = 0
$ws.Range("AC10").Value = 10
$ws.Range("AD10").Value = 6
$ws.Range("AK10").Value = 7
$ws.Range("AN10").Value = 16
$ws.Range("AO10").Value = 36
$ws.Range("AP10").Value = 1
$ws.Range("AQ10").Value = 53
$ws.Range("AR10").Value = $False
$ws.Range("AS10").Value = 10

# Row 11
$ws.Range("A11").Value = "FGFGHI314151"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 27
$ws.Range("F11").Value = 65
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = 83
$ws.Range("L11").Value = 99
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 6
$ws.Range("Q11").Value = 7
$ws.Range("R11").ClearContents()
$ws.Range("AC11").Value = 5
$ws.Range("AD11").Value = 8
$ws.Range("AI11").Value = 8
$ws.Range("AJ11").ClearContents()
$ws.Range("AN11").Value = 18
$ws.Range("AO11").Value = 36
$ws.Range("AP11").Value = 4
$ws.Range("AQ11").Value = 58
$ws.Range("AR11").Value = $False
$ws.Range("AS11").Value = 9
